# Lab01_ReviewReport.xlsx edit
#
# Fill in the reviewer name ("Onu Eduard Gabriel") and the review date
# (2019-03-18) in the merged D6:E6 / D7:E7 cells of all three
# "... Phase Defects" sheets, then leave the selection / active sheet
# the way it was left after doing the edit by hand:
#   - Requirements sheet ends up active, with D7:E7 selected
#   - Architect. Design sheet scrolled back to the top, E8 selected
#   - Coding sheet with H7 selected

$wb = $excel.ActiveWorkbook

$reviewerName = "Onu Eduard Gabriel"
$reviewDate   = Get-Date -Year 2019 -Month 3 -Day 18 -Hour 0 -Minute 0 -Second 0

$wsReq  = $wb.Worksheets.Item("Requirements Phase Defects")
$wsArch = $wb.Worksheets.Item("Architect. Design Phase Defects")
$wsCode = $wb.Worksheets.Item("Coding Phase Defects")

# --- Reviewer name (D6, shared across the three sheets) ---
$wsReq.Range("D6").Value  = $reviewerName
$wsArch.Range("D6").Value = $reviewerName
$wsCode.Range("D6").Value = $reviewerName

# --- Review date (D7) ---
# Set the number format first, then the value (setting the value first
# would let Excel auto-assign its own transient date format, left behind
# as an unused numFmt). Then copy the *format only* to the other two
# sheets so every sheet shares the same cellXfs entry instead of each
# NumberFormat assignment minting its own style.
$wsReq.Range("D7").NumberFormat = "mm-dd-yy"
$wsReq.Range("D7").Value = $reviewDate
$wsReq.Range("D7").Copy()
$wsArch.Range("D7").PasteSpecial(-4122)
$wsCode.Range("D7").PasteSpecial(-4122)

$wsArch.Range("D7").Value = $reviewDate
$wsCode.Range("D7").Value = $reviewDate

# --- Selection / active sheet state ---
$wsArch.Range("E8").Select()
$wsCode.Range("H7").Select()
$wsReq.Range("D7:E7").Select()
